$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Row 2 - reorder recorded-by list
$ws.Range("G2").Value = "gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, System, Veronia.rafat@med.asu.edu.eg"

# Row 3 - reorder recorded-by list
$ws.Range("G3").Value = "System, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"

# Row 4 - reorder recorded-by list
$ws.Range("G4").Value = "gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"

# Row 5 - new recorder added + updated student count
$ws.Range("G5").Value = "Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("H5").Value = "54/251"

# Row 7 - reorder recorded-by list
$ws.Range("G7").Value = "lamiaa.ossama@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg"

# Row 10 - updated average attendance % (use Formula + paste-values so the
# literal text is preserved instead of Excel auto-converting the "23.5%"
# string into a numeric percentage and reformatting the cell)
$c = $ws.Range("L10")
$c.Formula = "=""23.5%"""
$c.Copy()
$c.PasteSpecial(-4163)

# Row 15 - reorder recorded-by list + updated average attendance %
$ws.Range("G15").Value = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$c = $ws.Range("S15")
$c.Formula = "=""23.5%"""
$c.Copy()
$c.PasteSpecial(-4163)

# Row 28 - reorder recorded-by list
$ws.Range("G28").Value = "Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"

$excel.CutCopyMode = 0
